$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1288.6316
$ws.Range("I19").Value = 1348.5
$ws.Range("J19").Value = 1261
$ws.Range("K19").Value = 1348.5
$ws.Range("L19").Value = 1261
$ws.Range("M19").Value = -1173.5
$ws.Range("N19").Value = -1611
$ws.Range("H33").Value = 2088.8462
$ws.Range("I33").Value = 622.8570999999999
$ws.Range("J33").Value = 3799.1667
$ws.Range("K33").Value = 622.8570999999999
$ws.Range("L33").Value = 3799.1667
$ws.Range("M33").Value = -393.8570999999999
$ws.Range("N33").Value = -4257.1667
$ws.Range("H98").Value = 1021.2
$ws.Range("I98").Value = 1021.2
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1021.2
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 476.8
$ws.Range("N98").ClearContents()
$ws.Range("H107").Value = 279.2
$ws.Range("I107").Value = 234.5
$ws.Range("J107").Value = 458
$ws.Range("K107").Value = 234.5
$ws.Range("L107").Value = 458
$ws.Range("M107").Value = 1685.5
$ws.Range("N107").Value = -4298
$ws.Range("H122").Value = 1021.2
$ws.Range("I122").Value = 1021.2
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3063.6
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -613.6000000000004
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 6672407
$ws.Range("I6").Value = 12503875
$ws.Range("J6").Value = 7871.5713
$ws.Range("K6").Value = 12503875
$ws.Range("L6").Value = 7871.5713
$ws.Range("M6").Value = -12503702
$ws.Range("N6").Value = -8217.5713
$ws.Range("H10").Value = 2028.75
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 2028.75
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 2028.75
$ws.Range("N10").Value = -2368.75
$ws.Range("H19").Value = 1000
$ws.Range("I19").Value = 1100
$ws.Range("J19").Value = 900
$ws.Range("K19").Value = 1100
$ws.Range("L19").Value = 900
$ws.Range("M19").Value = -871
$ws.Range("N19").Value = -1358
$ws.Range("H26").Value = 4911.375
$ws.Range("I26").Value = 5158.2
$ws.Range("J26").Value = 4500
$ws.Range("K26").Value = 5158.2
$ws.Range("L26").Value = 4500
$ws.Range("M26").Value = -4828.2
$ws.Range("N26").Value = -5160
$ws.Range("H30").Value = 4333.3335
$ws.Range("I30").Value = 3000
$ws.Range("J30").Value = 5000
$ws.Range("K30").Value = 3000
$ws.Range("L30").Value = 5000
$ws.Range("M30").Value = -2850
$ws.Range("N30").Value = -5300
$ws.Range("H36").Value = 2500
$ws.Range("I36").Value = 2500
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 2500
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -2154
$ws.Range("H38").Value = 8903.799999999999
$ws.Range("I38").Value = 4839.6665
$ws.Range("J38").Value = 15000
$ws.Range("K38").Value = 4839.6665
$ws.Range("L38").Value = 15000
$ws.Range("M38").Value = -4372.6665
$ws.Range("N38").Value = -15934
$ws.Range("H39").Value = 7629
$ws.Range("I39").Value = 7629
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 7629
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -7109

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 1000000
$ws.Range("I7").Value = 1000000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1000000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -999887

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 1300
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 1300
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 1300
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -1640
$ws.Range("H17").Value = 4800
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 10500
$ws.Range("K17").Value = 1000
$ws.Range("L17").Value = 10500
$ws.Range("M17").Value = -826
$ws.Range("N17").Value = -10848
$ws.Range("H32").Value = 6333.3335
$ws.Range("I32").Value = 8500
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 8500
$ws.Range("L32").Value = 2000
$ws.Range("M32").Value = -8184
$ws.Range("N32").Value = -2632
$ws.Range("H35").Value = 1812.5
$ws.Range("I35").Value = 1812.5
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1812.5
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -1518.5
$ws.Range("H38").Value = 5574
$ws.Range("I38").Value = 4004.5
$ws.Range("J38").Value = 7666.6665
$ws.Range("K38").Value = 4004.5
$ws.Range("L38").Value = 7666.6665
$ws.Range("M38").Value = -3627.5
$ws.Range("N38").Value = -8420.666499999999
$ws.Range("H46").Value = 5574
$ws.Range("I46").Value = 4004.5
$ws.Range("J46").Value = 7666.6665
$ws.Range("K46").Value = 4004.5
$ws.Range("L46").Value = 7666.6665
$ws.Range("M46").Value = -3793.5
$ws.Range("N46").Value = -8088.6665
$ws.Range("H99").Value = 52497
$ws.Range("I99").Value = 4980
$ws.Range("J99").Value = 100014
$ws.Range("K99").Value = 4980
$ws.Range("L99").Value = 100014
$ws.Range("M99").Value = -3482
$ws.Range("N99").Value = -103010
$ws.Range("H126").Value = 52497
$ws.Range("I126").Value = 4980
$ws.Range("J126").Value = 100014
$ws.Range("K126").Value = 14940
$ws.Range("L126").Value = 300042
$ws.Range("M126").Value = -12470
$ws.Range("N126").Value = -304982

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 34.77778
$ws.Range("I12").Value = 5.5
$ws.Range("J12").Value = 38.4375
$ws.Range("K12").Value = 16.5
$ws.Range("L12").Value = 115.3125
$ws.Range("M12").Value = 156.5
$ws.Range("N12").Value = -461.3125
$ws.Range("H131").Value = 825.9184
$ws.Range("I131").Value = 522.63635
$ws.Range("J131").Value = 864.2643399999999
$ws.Range("K131").Value = 1567.90905
$ws.Range("L131").Value = 2592.79302
$ws.Range("M131").Value = 3472.09095
$ws.Range("N131").Value = -12672.79302

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 18250
$ws.Range("J6").Value = 18250
$ws.Range("L6").Value = 18250
$ws.Range("N6").Value = -18476
$ws.Range("H16").Value = 18250
$ws.Range("J16").Value = 18250
$ws.Range("L16").Value = 18250
$ws.Range("N16").Value = -18750
$ws.Range("H70").Value = 49568.49
$ws.Range("I70").Value = 78932.19
$ws.Range("J70").Value = 5522.9443
$ws.Range("K70").Value = 78932.19
$ws.Range("L70").Value = 5522.9443
$ws.Range("M70").Value = -78662.19
$ws.Range("N70").Value = -6062.9443
$ws.Range("H73").Value = 49568.49
$ws.Range("I73").Value = 78932.19
$ws.Range("J73").Value = 5522.9443
$ws.Range("K73").Value = 78932.19
$ws.Range("L73").Value = 5522.9443
$ws.Range("M73").Value = -77996.19
$ws.Range("N73").Value = -7394.9443

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 1800
$ws.Range("I32").Value = 1800
$ws.Range("K32").Value = 1800
$ws.Range("M32").Value = -1483
$ws.Range("H40").Value = 31512.94
$ws.Range("I40").Value = 50644
$ws.Range("K40").Value = 50644
$ws.Range("M40").Value = -50508
$ws.Range("H132").Value = 3045.3635
$ws.Range("I132").Value = 3372.1365
$ws.Range("J132").Value = 2391.818
$ws.Range("K132").Value = 10116.4095
$ws.Range("L132").Value = 7175.454000000001
$ws.Range("M132").Value = -7586.4095
$ws.Range("N132").Value = -12235.454

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 70011
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -70591
$ws.Range("H126").Value = 2392.6
$ws.Range("I126").Value = 2402.7646
$ws.Range("J126").Value = 2335
$ws.Range("K126").Value = 7208.293799999999
$ws.Range("L126").Value = 7005
$ws.Range("M126").Value = -4738.293799999999
$ws.Range("N126").Value = -11945
